# Commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
# Column G (header "K") holds per-outing strikeout counts. The values were
# regenerated (the sheet previously stored a different "Strike#" style
# figure in this column); write the newly computed K values for rows 2-70.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..70 (in row order)
$kValues = @(
    1,2,0,2,2,1,2,3,2,0,
    1,0,1,1,2,0,0,2,0,0,
    1,0,1,0,0,0,0,0,2,1,
    0,1,1,1,3,1,2,1,1,1,
    0,0,1,1,1,0,0,0,2,1,
    1,1,0,1,3,2,0,1,0,1,
    0,3,2,0,0,2,0,1,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
